$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "70.460.62"
Set-TextCell "E2" "  +1.61%  "
Set-TextCell "D3" "3.559.75"
Set-TextCell "E3" "  +2.08%  "
Set-TextCell "E4" "  +0.09%  "
Set-TextCell "D5" "588.02"
Set-TextCell "E5" "  +1.59%  "
Set-TextCell "D6" "188.95"
Set-TextCell "E6" "  +3.87%  "
Set-TextCell "D7" "3.551.42"
Set-TextCell "E7" "  +2.11%  "
Set-TextCell "E8" "  +1.98%  "
Set-TextCell "E9" "  -0.11%  "
Set-TextCell "E10" "  +10.87%  "
Set-TextCell "D11" "0.646"
Set-TextCell "E11" "  +1.13%  "
Set-TextCell "D12" "54.35"
Set-TextCell "E12" "  +1.50%  "
Set-TextCell "D13" "0.0000310"
Set-TextCell "E13" "  +3.13%  "
Set-TextCell "D14" "9.43"
Set-TextCell "E14" "  +0.85%  "
Set-TextCell "D15" "4.123.98"
Set-TextCell "E15" "  +2.36%  "
Set-TextCell "D16" "70.437.89"
Set-TextCell "E16" "  +1.76%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell "D17" "12.81"
Set-TextCell "E17" "  +4.80%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D18" "3.558.75"
Set-TextCell "E18" "  +1.75%  "
Set-TextCell "D19" "19.00"
Set-TextCell "E19" "  -0.72%  "
Set-TextCell "D20" "574.61"
Set-TextCell "E20" "  +7.42%  "
Set-TextCell "E21" "  +0.97%  "
Set-TextCell "D22" "0.996"
Set-TextCell "E22" "  -0.70%  "
Set-TextCell "D23" "18.00"
Set-TextCell "E23" "  -1.95%  "
Set-TextCell "D24" "4.64"
Set-TextCell "E24" "  +3.87%  "
Set-TextCell "E25" "  +0.54%  "
Set-TextCell "D26" "95.15"
Set-TextCell "E26" "  -0.71%  "
Set-TextCell "D27" "11.03"
Set-TextCell "E27" "  +0.34%  "
Set-TextCell "D28" "2.93"
Set-TextCell "E28" "  -0.15%  "
Set-TextCell "D29" "9.38"
Set-TextCell "E29" "  +3.99%  "
Set-TextCell "D30" "32.37"
Set-TextCell "E30" "  +1.54%  "
Set-TextCell "D31" "7.08"
Set-TextCell "E31" "  -1.05%  "
Set-TextCell "D32" "12.21"
Set-TextCell "E32" "  -1.51%  "
Set-TextCell "E33" "  +2.60%  "
Set-TextCell "D34" "3.82"
Set-TextCell "E34" "  +29.17%  "
Set-TextCell "D35" "63.35"
Set-TextCell "E35" "  -0.20%  "
Set-TextCell "D36" "3.23"
Set-TextCell "E36" "  +6.51%  "
$ws.Range("B37").Value = "TheGraph"
$ws.Range("C37").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextCell "D37" "0.408"
Set-TextCell "E37" "  +1.40%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D38" "527.39"
Set-TextCell "E38" "  -0.51%  "
Set-TextCell "D39" "3.692.47"
Set-TextCell "E39" "  +10.55%  "
Set-TextCell "D40" "38.11"
Set-TextCell "E40" "  +0.98%  "
Set-TextCell "E41" "  +0.04%  "
Set-TextCell "D42" "0.0₃0791"
Set-TextCell "E42" "  +5.57%  "
Set-TextCell "D43" "3.54"
Set-TextCell "E43" "  +5.80%  "
Set-TextCell "D44" "0.138"
Set-TextCell "E44" "  +3.54%  "
Set-TextCell "E45" "  +5.11%  "
Set-TextCell "D46" "3.47"
Set-TextCell "E46" "  -1.03%  "
Set-TextCell "E47" "  +0.49%  "
Set-TextCell "E48" "  +3.50%  "
Set-TextCell "D49" "9.23"
Set-TextCell "E49" "  +3.58%  "
Set-TextCell "D50" "0.999"
Set-TextCell "E50" "  +0.05%  "
Set-TextCell "E51" "  +7.99%  "
